# "Generate Report for Archive"
# Refresh the handoff-status report: cells that were still showing the old
# "Ready for handoff" status now reflect the current "In Translation" status,
# and the Status column on each sheet is narrowed to fit the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status columns E (zh-cn) and F (de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn / de-de detail sheets: Status column C
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Narrow the now-shorter Status columns to match the new content width
$newColumnWidth = 13.4101845877511

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
